# Natmi following Dr Hou advice
# Update the LR-pairs_lrc2p/Ccl24-Ccr2 sheet: target clusters now include
# FAPs and sCs in addition to ECs/M1/M2, and the underlying expression
# metrics are recomputed for the expanded set of clusters.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: M1 -> ECs
$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Ccl24"
$ws.Range("C2").Value = "Ccr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 0.8265333333333333
$ws.Range("H2").Value = 2.4796
$ws.Range("I2").Value = 0.2856046655980104
$ws.Range("J2").Value = 0.2856046655980104
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 0.1067053333333333
$ws.Range("N2").Value = 0.320116
$ws.Range("O2").Value = 0.0004703131782773351
$ws.Range("P2").Value = 0.000470313178277335
$ws.Range("Q2").Value = 0.08819551484444445
$ws.Range("R2").Value = 0.7937596336
$ws.Range("S2").Value = 0.0001343236380082358
$ws.Range("T2").Value = 0.0001343236380082357

# Row 3: M1 -> FAPs
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Ccl24"
$ws.Range("C3").Value = "Ccr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 0.8265333333333333
$ws.Range("H3").Value = 2.4796
$ws.Range("I3").Value = 0.2856046655980104
$ws.Range("J3").Value = 0.2856046655980104
$ws.Range("K3").Value = 2.0
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.08113566666666668
$ws.Range("N3").Value = 0.243407
$ws.Range("O3").Value = 0.0003576126147551242
$ws.Range("P3").Value = 0.0003576126147551241
$ws.Range("Q3").Value = 0.06706133302222222
$ws.Range("R3").Value = 0.6035519972000001
$ws.Range("S3").Value = 0.0001021358312507674
$ws.Range("T3").Value = 0.0001021358312507673

# Row 4: M1 -> M1
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Ccl24"
$ws.Range("C4").Value = "Ccr2"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 0.8265333333333333
$ws.Range("H4").Value = 2.4796
$ws.Range("I4").Value = 0.2856046655980104
$ws.Range("J4").Value = 0.2856046655980104
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 57.87945166666666
$ws.Range("N4").Value = 173.638355
$ws.Range("O4").Value = 0.2551087937213329
$ws.Range("P4").Value = 0.2551087937213328
$ws.Range("Q4").Value = 47.83929611755555
$ws.Range("R4").Value = 430.553665058
$ws.Range("S4").Value = 0.0728602617218931
$ws.Range("T4").Value = 0.07286026172189308

# Row 5: M1 -> M2
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Ccl24"
$ws.Range("C5").Value = "Ccr2"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 0.8265333333333333
$ws.Range("H5").Value = 2.4796
$ws.Range("I5").Value = 0.2856046655980104
$ws.Range("J5").Value = 0.2856046655980104
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 168.748281
$ws.Range("N5").Value = 506.244843
$ws.Range("O5").Value = 0.7437729482370158
$ws.Range("P5").Value = 0.7437729482370156
$ws.Range("Q5").Value = 139.4760791892
$ws.Range("R5").Value = 1255.2847127028
$ws.Range("S5").Value = 0.2124250241620792
$ws.Range("T5").Value = 0.2124250241620791

# Row 6: M1 -> sCs
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Ccl24"
$ws.Range("C6").Value = "Ccr2"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 0.8265333333333333
$ws.Range("H6").Value = 2.4796
$ws.Range("I6").Value = 0.2856046655980104
$ws.Range("J6").Value = 0.2856046655980104
$ws.Range("K6").Value = 2.0
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.065871
$ws.Range("N6").Value = 0.197613
$ws.Range("O6").Value = 0.0002903322486189976
$ws.Range("P6").Value = 0.0002903322486189975
$ws.Range("Q6").Value = 0.0544445772
$ws.Range("R6").Value = 0.4900011948
$ws.Range("S6").Value = 0.00008292024477914723
$ws.Range("T6").Value = 0.00008292024477914719

# Row 7: M2 -> ECs
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Ccl24"
$ws.Range("C7").Value = "Ccr2"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 2.067443666666667
$ws.Range("H7").Value = 6.202331
$ws.Range("I7").Value = 0.7143953344019897
$ws.Range("J7").Value = 0.7143953344019895
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 0.1067053333333333
$ws.Range("N7").Value = 0.320116
$ws.Range("O7").Value = 0.0004703131782773351
$ws.Range("P7").Value = 0.000470313178277335
$ws.Range("Q7").Value = 0.2206072655995556
$ws.Range("R7").Value = 1.985465390396
$ws.Range("S7").Value = 0.0003359895402690994
$ws.Range("T7").Value = 0.0003359895402690993

# Row 8: M2 -> FAPs
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Ccl24"
$ws.Range("C8").Value = "Ccr2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 2.067443666666667
$ws.Range("H8").Value = 6.202331
$ws.Range("I8").Value = 0.7143953344019897
$ws.Range("J8").Value = 0.7143953344019895
$ws.Range("K8").Value = 2.0
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.08113566666666668
$ws.Range("N8").Value = 0.243407
$ws.Range("O8").Value = 0.0003576126147551242
$ws.Range("P8").Value = 0.0003576126147551241
$ws.Range("Q8").Value = 0.1677434201907778
$ws.Range("R8").Value = 1.509690781717
$ws.Range("S8").Value = 0.0002554767835043569
$ws.Range("T8").Value = 0.0002554767835043567

# Row 9: M2 -> M1
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Ccl24"
$ws.Range("C9").Value = "Ccr2"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 2.067443666666667
$ws.Range("H9").Value = 6.202331
$ws.Range("I9").Value = 0.7143953344019897
$ws.Range("J9").Value = 0.7143953344019895
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 57.87945166666666
$ws.Range("N9").Value = 173.638355
$ws.Range("O9").Value = 0.2551087937213329
$ws.Range("P9").Value = 0.2551087937213328
$ws.Range("Q9").Value = 119.6625057783894
$ws.Range("R9").Value = 1076.962552005505
$ws.Range("S9").Value = 0.1822485319994398
$ws.Range("T9").Value = 0.1822485319994397

# Row 10: M2 -> M2
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Ccl24"
$ws.Range("C10").Value = "Ccr2"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 2.067443666666667
$ws.Range("H10").Value = 6.202331
$ws.Range("I10").Value = 0.7143953344019897
$ws.Range("J10").Value = 0.7143953344019895
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 168.748281
$ws.Range("N10").Value = 506.244843
$ws.Range("O10").Value = 0.7437729482370158
$ws.Range("P10").Value = 0.7437729482370156
$ws.Range("Q10").Value = 348.877564814337
$ws.Range("R10").Value = 3139.898083329033
$ws.Range("S10").Value = 0.5313479240749366
$ws.Range("T10").Value = 0.5313479240749364

# Row 11: M2 -> sCs
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Ccl24"
$ws.Range("C11").Value = "Ccr2"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 2.067443666666667
$ws.Range("H11").Value = 6.202331
$ws.Range("I11").Value = 0.7143953344019897
$ws.Range("J11").Value = 0.7143953344019895
$ws.Range("K11").Value = 2.0
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.065871
$ws.Range("N11").Value = 0.197613
$ws.Range("O11").Value = 0.0002903322486189976
$ws.Range("P11").Value = 0.0002903322486189975
$ws.Range("Q11").Value = 0.136184581767
$ws.Range("R11").Value = 1.225661235903
$ws.Range("S11").Value = 0.0002074120038398504
$ws.Range("T11").Value = 0.0002074120038398503

